# Requirements.xlsx update: clear a few obsolete requirement cells and
# refresh the sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: the "PID controller" comment is no longer applicable.
$ws.Range("C8").Value = ""

# Row 9: the two "toggle button / mode switching" comments are removed,
# and the row no longer needs the taller wrapped-text height.
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Rows.Item(9).RowHeight = 15

# Refresh the window scroll position / selection to match the author's
# final cursor position in the sheet.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C8").Select()
